# Moving from 2.0.1 to 2.0.2.
# Update the stack-trace text embedded in the document to reflect the new
# line numbers / object identity hash produced by the 2.0.2 run.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

# object identity hash changed
Replace-Text "EObjectImpl@34278424" "EObjectImpl@1cfde650"

# M2DocEvaluator.java line numbers changed
Replace-Text "M2DocEvaluator.caseQuery(M2DocEvaluator.java:555)" "M2DocEvaluator.caseQuery(M2DocEvaluator.java:559)"

# doSwitch appears three times at line 1096 -> 1216; wdReplaceAll (2) replaces
# every occurrence in the Content range in a single call.
Replace-Text "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1096)" "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1216)"

Replace-Text "M2DocEvaluator.caseBlock(M2DocEvaluator.java:1305)" "M2DocEvaluator.caseBlock(M2DocEvaluator.java:1425)"
Replace-Text "M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:283)" "M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:287)"
Replace-Text "M2DocEvaluator.generate(M2DocEvaluator.java:272)" "M2DocEvaluator.generate(M2DocEvaluator.java:276)"
